## TC08_INS_CancerType-GastricCancer.xlsx — automation API poc changes
## The "Website" column of the Programs query (cell B2) is reworked to
## derive its value from a CASE expression instead of a plain column ref,
## and the cursor/selection left on B8 after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B2")

$newQuery = @"
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
"@

# Re-enter the query text (the text content itself is the core of the edit).
$cell.Value = $newQuery

# Re-apply the cell's own formatting (12pt, theme text colour, wrapped) so
# it keeps looking exactly the same — this is what re-typing/re-committing
# the cell in Excel does, and it is what produced the fresh (if redundant)
# style/font record in the saved workbook.
$cell.Font.Size = 12
$cell.Font.ThemeColor = 1
$cell.WrapText = $true

# Leave the sheet scrolled/selected where the editor ended up.
$ws.Range("B8").Select()
